# Generate Report for Handback
#
# The handback run that produced "467e9683-a6ad-42cb-ae69-55bb3289f639.md"
# is superseded by a new run against "4a1395db-5321-4fac-8bad-e0cd38f24991.md"
# (updated timestamps / xliff hash), and a brand-new duplicate-content file
# "8be85514-958e-44cb-b01f-1ce8cc2cadd6.md" is handed back alongside it.
# This adds a second data row to each of the three report sheets
# (Overview / zh-cn / de-de) and refreshes the existing row's identifiers.

$wb = $excel.ActiveWorkbook

$oldGuid = "467e9683-a6ad-42cb-ae69-55bb3289f639"
$guid1   = "4a1395db-5321-4fac-8bad-e0cd38f24991"
$guid2   = "8be85514-958e-44cb-b01f-1ce8cc2cadd6"

$hash1 = "0fbeda4b754e85becf6f944f253f7cf126db9e56"
$hash2 = "30568b9836559814cf3a06dc968a6592c509bc49"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 2: refresh the file identifier + "Latest HO Xliff Generate Date"
$wsOv.Range("A2").Value = "$guid1.md"
$wsOv.Range("B2").Value = "e2e\$guid1.md"
$wsOv.Range("G2").Value = "2016-09-02 19:09:15"
foreach ($h in $wsOv.Hyperlinks) {
    if ($h.Range.Address() -eq "$`$B`$2") {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid1.md"
        $h.TextToDisplay = "e2e\$guid1.md"
    }
}

# Row 3 (new): add a ListRow to the Overview table and populate it
$loOv = $wsOv.ListObjects.Item("Overview")
$rowOv = $loOv.ListRows.Add()

$wsOv.Range("A3").Value = "$guid2.md"
$wsOv.Range("B3").Value = "e2e\$guid2.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-09-02 19:09:15"
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$hNew = $wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid2.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\$guid2.md")

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: refresh identifiers, xliff file name and the three timestamps
$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Range("G2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-02 19:09:09"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I2").Value = "$guid1.md"
$wsZh.Range("J2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-02 19:09:36"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "$`$A`$2") {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid1.md"
        $h.TextToDisplay = "$guid1.md"
    } elseif ($addr -eq "$`$I`$2") {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/948151c2f074c468f759ff0b3450616702a29cd6/e2e/$guid1.md"
        $h.TextToDisplay = "$guid1.md"
    }
}

# Row 3 (new): add a ListRow to the zh-cn table and populate it
$loZh = $wsZh.ListObjects.Item("zh-cn")
$rowZh = $loZh.ListRows.Add()

$wsZh.Range("A3").Value = "$guid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-02 19:09:09"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = "$guid2.md"
$wsZh.Range("J3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-02 19:09:36"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$hZhA = $wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid2.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$guid2.md")
$hZhI = $wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/948151c2f074c468f759ff0b3450616702a29cd6/e2e/$guid2.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$guid2.md")

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: refresh identifiers, xliff file name and the three timestamps
$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Range("G2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-02 19:09:15"
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I2").Value = "$guid1.md"
$wsDe.Range("J2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-02 19:09:44"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "$`$A`$2") {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid1.md"
        $h.TextToDisplay = "$guid1.md"
    } elseif ($addr -eq "$`$I`$2") {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c83e3ab1a0f2747369d9c591bec928fc2413a958/e2e/$guid1.md"
        $h.TextToDisplay = "$guid1.md"
    }
}

# Row 3 (new): add a ListRow to the de-de table and populate it
$loDe = $wsDe.ListObjects.Item("de-de")
$rowDe = $loDe.ListRows.Add()

$wsDe.Range("A3").Value = "$guid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-02 19:09:15"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = "$guid2.md"
$wsDe.Range("J3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-02 19:09:44"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$hDeA = $wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682f8c696bcb3356602fed37c6c5c9f317e9be55/e2e/$guid2.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$guid2.md")
$hDeI = $wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c83e3ab1a0f2747369d9c591bec928fc2413a958/e2e/$guid2.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$guid2.md")

Write-Output "done"
